# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 1.78 = 6518.26 pesos
✅ 6518.26 pesos = 1.77 = 945.96 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$newText = $newText.TrimEnd("`r", "`n")
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 561.2
$wsTasas.Range("O10").Value = 3658.05
$wsTasas.Range("N12").Value = 3678.9
$wsTasas.Range("O12").Value = 533.9
